$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns
# D-column price cells are forced to Text format to preserve the
# literal "dotted" formatting (e.g. "37.824.66") used by this sheet,
# matching the original inlineStr text values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.824.66'
$ws.Range("E2").Value = '  +0.06%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.090.66'
$ws.Range("E3").Value = '  +0.16%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.94'
$ws.Range("E5").Value = '  -0.30%  '

$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.63'
$ws.Range("E7").Value = '  -0.37%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.395'
$ws.Range("E9").Value = '  +0.49%  '

$ws.Range("E11").Value = '  +2.95%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.13'
$ws.Range("E12").Value = '  +2.40%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.398.09'
$ws.Range("E13").Value = '  +0.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.41'
$ws.Range("E14").Value = '  +0.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.782'
$ws.Range("E15").Value = '  +1.45%  '

$ws.Range("E16").Value = '  +0.97%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.089.66'
$ws.Range("E17").Value = '  +0.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.825.85'
$ws.Range("E18").Value = '  +0.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.15'
$ws.Range("E19").Value = '  -0.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.30'
$ws.Range("E20").Value = '  -0.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0839'
$ws.Range("E21").Value = '  -0.05%  '

$ws.Range("E22").Value = '  +0.77%  '

$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.40'
$ws.Range("E24").Value = '  -0.36%  '

$ws.Range("E25").Value = '  +0.73%  '

$ws.Range("E26").Value = '  +9.34%  '

$ws.Range("E27").Value = '  +0.98%  '

$ws.Range("E28").Value = '  -2.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.55'
$ws.Range("E29").Value = '  +0.05%  '

$ws.Range("E30").Value = '  -0.87%  '

$ws.Range("E31").Value = '  +1.07%  '

$ws.Range("E32").Value = '  +0.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.68'
$ws.Range("E34").Value = '  -0.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.50'
$ws.Range("E35").Value = '  -0.50%  '

$ws.Range("E36").Value = '  -0.56%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.41'
$ws.Range("E37").Value = '  -1.54%  '

$ws.Range("E38").Value = '  +0.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.38'
$ws.Range("E39").Value = '  -0.18%  '

$ws.Range("E40").Value = '  +10.74%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '102.73'
$ws.Range("E41").Value = '  +3.88%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0977'
$ws.Range("E42").Value = '  -1.80%  '

$ws.Range("E43").Value = '  -0.52%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.86'
$ws.Range("E44").Value = '  +5.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.458.10'
$ws.Range("E45").Value = '  -0.57%  '

$ws.Range("E46").Value = '  -0.84%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.19'
$ws.Range("E47").Value = '  -4.53%  '

$ws.Range("E48").Value = '  -0.68%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.33'
$ws.Range("E49").Value = '  -1.10%  '

$ws.Range("E50").Value = '  -1.56%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.282.77'
$ws.Range("E51").Value = '  +0.18%  '
